# Regenerate merged AHB files
# 1. Rename the "_old" / "_new" header suffixes to "_FV2304" / "_FV2310"
# 2. Turn the data range into an Excel Table ("Table1")
# 3. Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (shared-string text) -------------------------
$oldCols = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldCols.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($oldCols[$i] + "_FV2304")
}

for ($i = 0; $i -lt $oldCols.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($oldCols[$i] + "_FV2310")
}

# --- 2. Convert the range A1:U58 into a real table ------------------------
$range = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze panes at row 1 ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
